$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "E2"  = "  -2.34%  "
    "E3"  = "  -5.10%  "
    "E4"  = "  -0.04%  "
    "E6"  = "  -7.59%  "
    "E7"  = "  +0.00%  "
    "E8"  = "  -10.73%  "
    "E9"  = "  -5.05%  "
    "E10" = "  -2.19%  "
    "E11" = "  -0.12%  "
    "E12" = "  -2.52%  "
    "E13" = "  -3.61%  "
    "E14" = "  -6.85%  "
    "E15" = "  -5.07%  "
    "E16" = "  -1.60%  "
    "E17" = "  -4.55%  "
    "E18" = "  -5.10%  "
    "E19" = "  -4.14%  "
    "E20" = "  -1.30%  "
    "E21" = "  -2.67%  "
    "E22" = "  -6.89%  "
    "E23" = "  -0.14%  "
    "E24" = "  -2.75%  "
    "E25" = "  -1.58%  "
    "E26" = "  +10.26%  "
    "E27" = "  +0.30%  "
    "E28" = "  -5.05%  "
    "E29" = "  -9.52%  "
    "E30" = "  -4.06%  "
    "E31" = "  -6.56%  "
    "E32" = "  -6.55%  "
    "E33" = "  -2.09%  "
    "E34" = "  -5.02%  "
    "E35" = "  -5.69%  "
    "E36" = "  -0.01%  "
    "E37" = "  -3.94%  "
    "E38" = "  -2.17%  "
    "E39" = "  -0.08%  "
    "E40" = "  -9.36%  "
    "E41" = "  +2.38%  "
    "E42" = "  -0.11%  "
    "E43" = "  -1.72%  "
    "E44" = "  -0.88%  "
    "E45" = "  -1.97%  "
    "E46" = "  -7.64%  "
    "E47" = "  -2.62%  "
    "E48" = "  -4.58%  "
    "E49" = "  -11.21%  "
    "E50" = "  -4.15%  "
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
